$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize double spaces to single spaces and "and" -> "+" before "listening"
# in the monthly comment column (C), per the Feb-report copy edit.

$ws.Range("C2").Value2 = 'Well done this month Ai! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Unicorn!”). You’re really confident now when speaking English, and your phonics + listening is getting very good!'
$ws.Range("C3").Value2 = 'Well done this month Aoi! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re always so polite & well behaved, and your phonics + listening is getting very good!'
$ws.Range("C5").Value2 = 'Excellent work this month Ayumu! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Knight!”). You’ve been so friendly & get along well with the other kids in class, and your phonics + listening is getting very good!'
$ws.Range("C6").Value2 = 'Great job this month Beni! You’re always singing and answering with a big voice during calendar time. You’ve been so friendly & get along well with the other kids in class, and your phonics + listening is getting very good!'
$ws.Range("C7").Value2 = 'Excellent work this month Emi! You’re always singing and answering with a big voice during calendar time. You’ve been so friendly & get along well with the other kids in class, and your phonics + listening is getting very good!'
$ws.Range("C9").Value2 = 'Well done this month Hiroto! It’s great to hear you using English in class during snack time, and when I hand your workbooks back. I saw you really enjoyed the Pinocchio song & dance, you were having so much fun doing the silly actions – great job!'
$ws.Range("C10").Value2 = 'Excellent work this month Hisui! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re really confident now when speaking English, and your phonics + listening is getting very good!'
$ws.Range("C11").Value2 = 'Excellent work this month Hiyori! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re really confident now when speaking English, and your phonics + listening is getting very good!'
$ws.Range("C14").Value2 = 'Well done this month Kansuke! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Knight!”). You’re really confident when speaking English, and you did a great job both writing and reading the 3 letter words on the handwriting worksheet, excellent!'
$ws.Range("C15").Value2 = 'Well done this month Kanta! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Knight!”). You’re really confident when speaking English, and you did a great job both writing and reading the 3 letter words on the handwriting worksheet, excellent!'
$ws.Range("C17").Value2 = 'Awesome work this month Kimueru! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re always so polite & well behaved, and your phonics + listening is getting very good!'
$ws.Range("C20").Value2 = 'Excellent work this month Manaka! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re really confident now when speaking English, and your phonics + listening is getting very good!'
$ws.Range("C23").Value2 = 'Excellent work this month Mengtao! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re always so polite & well behaved, and your phonics + listening is getting very good!'
$ws.Range("C24").Value2 = 'Excellent work this month Mio! You had so much fun singing & dancing to the Pinocchio song! You’re really confident now when speaking English, and you did a great job remembering the Valentines Day words & phrases!'
$ws.Range("C25").Value2 = 'Well done this month Mio! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re really confident now when speaking English, Tomomi & I are so impressed with your big voice, especially during calendar time!'
$ws.Range("C26").Value2 = 'Excellent work this month Miran! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re always so polite & well behaved, and your phonics + listening is getting very good!'
$ws.Range("C28").Value2 = 'Awesome work this month Reika! You had so much fun singing & dancing to the Pinocchio song! You’re always so polite & well behaved, and you did a great job remembering the Valentines Day words & phrases!'
$ws.Range("C32").Value2 = 'Well done this month Ryosei! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Knight!”). You’re always so friendly & get along well with the other kids in class, and you learned so many of the fairy tale flash cards, you could even help me describe many of them!'
$ws.Range("C34").Value2 = 'Well done this month Saki! You had so much fun singing & dancing to the Pinocchio song! You’re always so polite & well behaved, and you learned so many of the fairy tale flash cards, you could even help me describe some of them!'
$ws.Range("C35").Value2 = 'Excellent work this month Sara! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Unicorn!”). You’re really confident now when speaking English, and you learned most of the fairy tale flash cards, you could even help me describe some of them!'
$ws.Range("C36").Value2 = 'Excellent work this month Shotaro! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re really confident now when speaking English, and your phonics + listening is getting very good!'
$ws.Range("C38").Value2 = 'Awesome work this month Soma! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re always so polite & well behaved, and your phonics + listening is getting very good!'
$ws.Range("C39").Value2 = 'Well done this month Taichi! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Knight!”). You’ve been so friendly & get along well with the other kids in class, and you learned so many of the fairy tale flash cards, you could even help me describe some of them!'
$ws.Range("C45").Value2 = 'Excellent work this month Yui! It’s so great to see your writing and reading improve, you did a great job on the handwriting sheet this month! You’re always so polite & well behaved, and your phonics + listening is getting very good!'
$ws.Range("C51").Value2 = 'Excellent work this month Yuuri! You were so fast to pick up the Q&A this month (“What do you want to be?” “I want to be a Knight!”). You’ve been so friendly & get along well with the other kids in class, and you learned so many of the fairy tale flash cards, you could even help me describe some of them!'

# Update the saved selection/active cell to match the editor's last position
$ws.Range("C32").Select()

